{"js": "// Bug fix: \"pour\" and \"mouler\" were split across three runs (\"pour \", \" \",\n// \"mouler\") by mk_frq_cnts, leaving a stray extra space and broken run\n// boundaries. Re-join them into a single text node \"pour mouler\" so the\n// word is represented as one run (using the formatting of the original\n// \"pour \" run).\nconst body = context.document.body;\n\n// The buggy split always leaves a double space between \"pour\" and\n// \"mouler\" (the lone middle run containing just \" \"). Search for that\n// exact, unique pattern so we only touch the affected text.\nconst results = body.search(\"pour  mouler\", { matchCase: true, matchWildcards: false });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find the split 'pour  mouler' text to fix.\");\n}\n\nfor (let i = 0; i < results.items.length; i++) {\n  const found = results.items[i];\n  found.insertText(\"pour mouler\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Bug fix: \"pour\" and \"mouler\" were split across three runs (\"pour \", \" \",\n# \"mouler\") by mk_frq_cnts, leaving a stray lone-space run between them.\n# Re-join them into a single text node \"pour mouler\" (one run, carrying the\n# formatting of the original \"pour \" run).\n\n$d = $word.ActiveDocument\n\n# The buggy split always leaves a double space between \"pour\" and \"mouler\"\n# (the lone middle run contains just that extra \" \"). This exact phrase is\n# unique in the document, so Find reliably locates the affected text.\n$find = $d.Content.Find\n$find.Text = \"pour  mouler\"\n$find.MatchCase = $true\n$find.MatchWildcards = $false\n$found = $find.Execute()\n\nif (-not $found) {\n    throw \"Could not find the split 'pour  mouler' text to fix.\"\n}\n\n$matchRange = $find.Parent\n$matchStart = $matchRange.Start\n$matchEnd = $matchRange.End\n\n# Locate \"pour\" within the match so we know exactly where its own (legitimate)\n# trailing space ends, independent of any hard-coded offsets.\n$pourFind = $d.Range($matchStart, $matchEnd).Find\n$pourFind.Text = \"pour\"\n$pourFind.MatchCase = $true\n$pourFind.Execute() | Out-Null\n$afterPour = $pourFind.Parent.End\n\n# Right after \"pour\" there are two consecutive spaces in the flattened text:\n# index 0 is the trailing space that belongs to the \"pour \" run, index 1 is\n# the stray extra space living in its own separate run. Deleting just that\n# stray character merges the surrounding \"pour \" and \"mouler\" runs into a\n# single \"pour mouler\" run (same formatting as the original \"pour \" run).\n$strayRange = $d.Range($afterPour + 1, $afterPour + 2)\n$strayRange.Delete()\n"}
